$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values for rows 2, 3, 4, and 9
$ws.Range("F2").Value = 7
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 2
$ws.Range("F9").Value = -9
